$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price-report record needs to be inserted as row 115 ("Early Burlat" /
# "Segunda"), pushing all the existing rows from 115..183 down to 116..184.
$ws.Rows.Item(115).Insert()

# Fill in the newly inserted row with the new record's data (columns A-J are
# identical to the rest of this Macroferia Regional de Talca / Cereza block).
$ws.Range("A115").Value = 5
$ws.Range("B115").Value = "Macroferia Regional de Talca"
$ws.Range("C115").Value = "Maule"
$ws.Range("D115").Value = 44873
$ws.Range("E115").Value = 7
$ws.Range("F115").Value = "Fruta"
$ws.Range("G115").Value = 100103
$ws.Range("H115").Value = "Frutos de hueso (carozo)"
$ws.Range("I115").Value = 100103001
$ws.Range("J115").Value = "Cereza"
$ws.Range("K115").Value = "Early Burlat"
$ws.Range("L115").Value = "Segunda"
$ws.Range("M115").Value = 180
$ws.Range("N115").Value = 19000
$ws.Range("O115").Value = 20000
$ws.Range("P115").Value = 19667
$ws.Range("Q115").Value = "$/bandeja 10 kilos"
$ws.Range("R115").Value = "Provincia de Curicó"
$ws.Range("S115").Value = 1967
$ws.Range("T115").Value = 10
